# "replaced all reading times with average"
#
# For each of the three per-model sheets (10_trees, 500_trees, 1600_trees):
#   - column K ("total framework time") used to be SUM(E,G,H) for that row;
#     it now reads SUM($E$47,G,H) -- i.e. it always reads the *average*
#     "data loading time" (a new row 47) instead of each row's own value.
#   - a new row 47 is appended: D47 = "average", E47 = the average data
#     loading time. On 10_trees this is the real AVERAGE(...) formula;
#     on 500_trees / 1600_trees it's just a link back to '10_trees'!E47
#     (so all three sheets agree on the same average).
#   - 500_trees had a stray, nearly-empty row 46 (K46 only) that is removed.
#
# L (=K/1000) is unaffected in formula text; it just recalculates because K
# changed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("10_trees", "500_trees", "1600_trees")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # 500_trees has a leftover row 46 (only K46 = SUM(E46,G46,H46) = 0)
    # that needs to disappear before we lay down the new row 47.
    if ($name -eq "500_trees") {
        $ws.Rows.Item(46).Delete()
    }

    # Row 2's K formula is written out (not shared); K3:K45 is one shared
    # formula. Re-assigning the whole range re-derives a clean shared
    # formula with the same relative G/H but an absolute $E$47.
    $ws.Cells.Item(2, 11).Formula = "=SUM(`$E`$47,G2,H2)"
    $ws.Range("K3:K45").Formula = "=SUM(`$E`$47,G3,H3)"

    # New summary row.
    $ws.Cells.Item(47, 4).Value = "average"
}

$ws1 = $wb.Worksheets.Item("10_trees")
$ws1.Cells.Item(47, 5).Formula = "=AVERAGE(E2:E45,'500_trees'!E2:E45,'10_trees'!E2:E45)"

$ws2 = $wb.Worksheets.Item("500_trees")
$ws2.Cells.Item(47, 5).Formula = "='10_trees'!E47"

$ws3 = $wb.Worksheets.Item("1600_trees")
$ws3.Cells.Item(47, 5).Formula = "='10_trees'!E47"

$excel.CalculateFull()

# Restore selections / scroll position roughly as in the target file, and
# finish with "final" active again (it was, and stays, the active tab).
$ws1.Activate()
$ws1.Range("K2:K45").Select()

$ws2.Activate()
$ws2.Range("O43").Select()

$ws3.Activate()
$ws3.Range("N39").Select()

$wsf = $wb.Worksheets.Item("final")
$wsf.Activate()
$wsf.Range("I29").Select()
